$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G column values for rows 15-19, 24-25
$ws.Range("G15").Value = 58
$ws.Range("G16").Value = 54
$ws.Range("G17").Value = 54.5
$ws.Range("G18").Value = 58
$ws.Range("G19").Value = 58
$ws.Range("G24").Value = 58
$ws.Range("G25").Value = 55.2

# Update the active cell selection to G26
$ws.Range("G26").Select()
